# "Improve plot of CrankFitting"
# Adds four new species rows (Hexanal, 2-methylbutanal, 2-methylpropanal,
# Pyridine) with their permeation/diffusion data and source-paper
# citations, replaces the hard-coded formula in D3 with its resulting
# literal value, highlights two "total concentration" cells that are
# estimates, and touches up a couple of cosmetic view/format details.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- D3: bake the formula result into a literal value ------------------
$ws.Range("D3").Value = 0.139

# --- Row 5: Hexanal ------------------------------------------------------
$ws.Range("A5").Value = "Hexanal"
$ws.Range("B5").Value = [double]"2.796E-10"
$ws.Range("C5").Value = [double]"5.5570000000000001E-8"
$ws.Range("C5").NumberFormat = "0.00E+00"

# --- Row 6: 2-methylbutanal ----------------------------------------------
$ws.Range("A6").Value = "2-methylbutanal"
$ws.Range("B6").Value = [double]"5.5799999999999997E-8"
$ws.Range("C6").Value = [double]"4.7027000000000001E-8"
$ws.Range("C6").NumberFormat = "0.00E+00"

# --- Row 7: 2-methylpropanal ----------------------------------------------
$ws.Range("A7").Value = "2-methylpropanal"
$ws.Range("B7").Value = [double]"2.4E-8"
$ws.Range("C7").Value = [double]"2.5424E-8"
$ws.Range("C7").NumberFormat = "0.00E+00"
$ws.Range("J7").Value = "Investigation of roasted coffee freshness with an improved headspace technique"

# --- Row 5 (paper ref + alignment reset) ----------------------------------
$ws.Range("J5").Value = "Coffee roasting and quenching technology -formation and stability of aroma compounds"
$ws.Range("J5").WrapText = $False

# --- Row 5 (highlighted, estimated concentration) + rest of row 5 --------
$ws.Range("D5").Value = 0.01
$ws.Range("D5").Interior.Color = 65535
$ws.Range("E5").Value = 10
$ws.Range("F5").Value = [double]"5.1900000000000004E-4"
$ws.Range("F5").NumberFormat = "0.00E+00"
$ws.Range("G5").Value = 100.161
$ws.Range("G5").NumberFormat = "0.00E+00"
$ws.Range("H5").Value = 0
$ws.Range("H5").NumberFormat = "0.00E+00"
$ws.Range("I5").Value = 0
$ws.Range("I5").NumberFormat = "0.00E+00"

# --- Rest of row 6 ---------------------------------------------------------
$ws.Range("D6").Value = [double]"2.07E-2"
$ws.Range("E6").Value = [double]"49.317000999999998"
$ws.Range("F6").Value = [double]"4.9346000000000004E-4"
$ws.Range("F6").NumberFormat = "0.00E+00"
$ws.Range("G6").Value = 86.134
$ws.Range("G6").NumberFormat = "0.00E+00"
$ws.Range("H6").Value = 0
$ws.Range("H6").NumberFormat = "0.00E+00"
$ws.Range("I6").Value = 0
$ws.Range("I6").NumberFormat = "0.00E+00"
$ws.Range("J6").Value = "Coffee roasting and quenching technology -formation and stability of aroma compounds"
$ws.Range("J6").WrapText = $False
$ws.Rows.Item(6).RowHeight = 16.2

# --- Rest of row 7 ---------------------------------------------------------
$ws.Range("D7").Value = [double]"1.7399999999999999E-2"
$ws.Range("E7").Value = 170
$ws.Range("F7").Value = [double]"1.9699999999999999E-4"
$ws.Range("F7").NumberFormat = "0.00E+00"
$ws.Range("G7").Value = [double]"72.058000000000007"
$ws.Range("G7").NumberFormat = "0.00E+00"
$ws.Range("H7").Value = 0
$ws.Range("H7").NumberFormat = "0.00E+00"
$ws.Range("I7").Value = 0
$ws.Range("I7").NumberFormat = "0.00E+00"

# --- Row 8: Pyridine --------------------------------------------------------
$ws.Range("A8").Value = "Pyridine"
$ws.Range("B8").Value = [double]"4.8208374409751881E-8"
$ws.Range("C8").Value = [double]"8.7099999999999998E-12"
$ws.Range("C8").NumberFormat = "0.00E+00"
$ws.Range("D8").Value = [double]"2.3E-3"
$ws.Range("D8").Interior.Color = 65535
$ws.Range("E8").Value = 20.8
$ws.Range("F8").Value = [double]"1.1E-5"
$ws.Range("F8").NumberFormat = "0.00E+00"
$ws.Range("G8").Value = [double]"79.102000000000004"
$ws.Range("G8").NumberFormat = "0.00E+00"
$ws.Range("H8").Value = 0
$ws.Range("H8").NumberFormat = "0.00E+00"
$ws.Range("I8").Value = 0
$ws.Range("I8").NumberFormat = "0.00E+00"
$ws.Range("J8").Value = "Approximate permeation as benzene"

# --- Cosmetic touch-ups ------------------------------------------------
$ws.Rows.Item(1).RowHeight = 96.9
$ws.Columns.Item(1).ColumnWidth = 15.14
[void]$ws.Range("K11").Select()

Write-Host "Applied CrankFitting data update"
